# Update cryptos list with latest price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $ref, $val) {
    # Force Excel to store the value as text instead of silently
    # re-parsing it as a number (which would drop the original
    # text formatting, e.g. trailing zeros or thousand-dot
    # grouping), then restore the default "Normal" style so the
    # cell keeps its original (unstyled) appearance.
    $c = $ws.Range($ref)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue $ws "D2" '58.378.48'
$ws.Range("E2").Value = '  -2.92%  '
Set-TextValue $ws "D3" '2.280.16'
$ws.Range("E3").Value = '  -5.68%  '
$ws.Range("E4").Value = '  -0.03%  '
Set-TextValue $ws "D5" '545.21'
$ws.Range("E5").Value = '  -1.36%  '
Set-TextValue $ws "D6" '130.71'
$ws.Range("E6").Value = '  -4.68%  '
$ws.Range("E7").Value = '  -0.01%  '
Set-TextValue $ws "D8" '0.570'
$ws.Range("E8").Value = '  -2.85%  '
Set-TextValue $ws "D9" '2.278.61'
$ws.Range("E9").Value = '  -5.60%  '
$ws.Range("E10").Value = '  -3.83%  '
Set-TextValue $ws "D11" '5.50'
$ws.Range("E11").Value = '  -2.79%  '
$ws.Range("E12").Value = '  +0.07%  '
$ws.Range("E13").Value = '  -5.37%  '
Set-TextValue $ws "D14" '23.57'
$ws.Range("E14").Value = '  -5.31%  '
Set-TextValue $ws "D15" '2.687.26'
$ws.Range("E15").Value = '  -5.72%  '
Set-TextValue $ws "D16" '58.362.01'
$ws.Range("E16").Value = '  -2.76%  '
$ws.Range("E17").Value = '  -3.42%  '
Set-TextValue $ws "D18" '2.279.11'
$ws.Range("E18").Value = '  -5.56%  '
Set-TextValue $ws "D19" '10.61'
$ws.Range("E19").Value = '  -6.08%  '
Set-TextValue $ws "D20" '4.29'
$ws.Range("E20").Value = '  -4.11%  '
Set-TextValue $ws "D21" '314.22'
$ws.Range("E21").Value = '  -4.07%  '
$ws.Range("E22").Value = '  -4.54%  '
$ws.Range("E23").Value = '  +0.12%  '
Set-TextValue $ws "D24" '62.91'
$ws.Range("E24").Value = '  -3.69%  '
$ws.Range("E25").Value = '  -3.93%  '
Set-TextValue $ws "D26" '1.00'
$ws.Range("E26").Value = '  -0.42%  '
Set-TextValue $ws "D27" '8.10'
$ws.Range("E27").Value = '  -6.80%  '
$ws.Range("E28").Value = '  -6.55%  '
$ws.Range("E29").Value = '  -0.61%  '
Set-TextValue $ws "D30" '170.58'
$ws.Range("E30").Value = '  +0.33%  '
Set-TextValue $ws "D31" '0.0₃0722'
$ws.Range("E31").Value = '  -6.46%  '
$ws.Range("E32").Value = '  -0.39%  '
Set-TextValue $ws "D33" '5.75'
$ws.Range("E33").Value = '  -5.66%  '
Set-TextValue $ws "D34" '0.384'
$ws.Range("E34").Value = '  -4.92%  '
$ws.Range("E36").Value = '  -4.03%  '
$ws.Range("E37").Value = '  +0.06%  '
$ws.Range("E38").Value = '  -5.00%  '
$ws.Range("E39").Value = '  -5.93%  '
Set-TextValue $ws "D40" '37.89'
$ws.Range("E40").Value = '  -2.33%  '
Set-TextValue $ws "D41" '1.51'
$ws.Range("E41").Value = '  -5.59%  '
Set-TextValue $ws "D42" '294.50'
$ws.Range("E42").Value = '  -10.47%  '
Set-TextValue $ws "D43" '140.80'
$ws.Range("E43").Value = '  -3.16%  '
$ws.Range("E44").Value = '  -5.73%  '
Set-TextValue $ws "D45" '0.0947'
$ws.Range("E45").Value = '  -1.76%  '
$ws.Range("E46").Value = '  -3.42%  '
Set-TextValue $ws "D47" '0.553'
$ws.Range("E47").Value = '  -3.87%  '
Set-TextValue $ws "D48" '18.34'
$ws.Range("E48").Value = '  -8.35%  '
$ws.Range("E49").Value = '  -4.23%  '
Set-TextValue $ws "D50" '16.59'
$ws.Range("E50").Value = '  -5.43%  '
Set-TextValue $ws "D51" '10.99'
$ws.Range("E51").Value = '  -0.48%  '
